$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "System Information" heading: center it and bump font size 12pt -> 16pt
#    (w:sz / w:szCs 24 -> 32), keeping bold + single underline.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "System Information") {
        $p.Alignment = 1          # wdAlignParagraphCenter
        $p.Range.Font.Size = 16   # half-points 24 -> 32
        $p.Range.Font.SizeBi = 16 # also updates szCs
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Merge the split "2 inner rocky planets and thre" + bookmark + "e outer
#    gas giants orbit " runs into a single contiguous run/text (this also
#    drops the stray mid-word _GoBack bookmark that used to live there).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "2 inner rocky planets and three outer gas giants orbit ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2 inner rocky planets and three outer gas giants orbit ", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 3) After "Rains slightly more than on earth." add two new bullet
#    paragraphs:
#      a) ilvl=2 bullet with text "There are many different forms of life
#         that inhabit the planet"
#      b) ilvl=3 (empty) bullet that now hosts the relocated _GoBack
#         bookmark.
# ---------------------------------------------------------------------------
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text.Trim() -eq "Rains slightly more than on earth.") {
        $targetIdx = $idx
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$target.Range.InsertParagraphAfter() | Out-Null

# -- new paragraph #1 (ilvl=2) -------------------------------------------------
$newPara1 = $d.Paragraphs.Item($targetIdx + 1)
$newPara1.Range.ListFormat.ListLevelNumber = 3   # 1-based -> w:ilvl=2
$newPara1.Range.Text = "There are many different forms of life that inhabit the planet"

# -- new paragraph #2 (ilvl=3, empty, carries the _GoBack bookmark) -----------
$newPara1 = $d.Paragraphs.Item($targetIdx + 1)
$newPara1.Range.InsertParagraphAfter() | Out-Null

$newPara2 = $d.Paragraphs.Item($targetIdx + 2)
$newPara2.Range.ListFormat.ListLevelNumber = 4   # 1-based -> w:ilvl=3

# The bookmark engine in this host only emits a matching bookmarkEnd for a
# non-collapsed range, so stamp a one-character placeholder, bookmark it,
# then delete just that character (leaving start/end adjacent & the
# paragraph empty, matching the target markup).
$newPara2.Range.Text = "X"
$newPara2 = $d.Paragraphs.Item($targetIdx + 2)
$placeholder = $d.Range($newPara2.Range.Start, $newPara2.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null

$newPara2 = $d.Paragraphs.Item($targetIdx + 2)
$charRange = $d.Range($newPara2.Range.Start, $newPara2.Range.Start + 1)
$charRange.Delete() | Out-Null
